$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data for rows 2-19 (columns A=Player, B=Position, C=Team),
# reflecting the reshuffled roster order from the commit.
$data = @(
    @("Jalen Brunson",      "PG",       "New York Knicks"),
    @("Devin Booker",       "PG,SG",    "Phoenix Suns"),
    @("Devin Vassell",      "SG,SF",    "San Antonio Spurs"),
    @("Desmond Bane",       "SG,SF",    "Memphis Grizzlies"),
    @("Jeremy Sochan",      "SF,PF",    "San Antonio Spurs"),
    @("P.J. Washington",    "PF",       "Dallas Mavericks"),
    @("Alperen Sengün",     "C",        "Houston Rockets"),
    @("Dereck Lively II",   "C",        "Dallas Mavericks"),
    @("Trae Young",         "PG",       "Atlanta Hawks"),
    @("Jared McCain",       "PG,SG",    "Philadelphia 76ers"),
    @("LeBron James",       "SF,PF",    "Los Angeles Lakers"),
    @("Nicolas Claxton",    "C",        "Brooklyn Nets"),
    @("Coby White",         "PG,SG",    "Chicago Bulls"),
    @("Walker Kessler",     "C",        "Utah Jazz"),
    @("Shaedon Sharpe",     "SG,SF",    "Portland Trail Blazers"),
    @("Immanuel Quickley",  "PG,SG",    "Toronto Raptors"),
    @("Kawhi Leonard",      "SG,SF,PF", "LA Clippers"),
    @("Norman Powell",      "SG,SF",    "LA Clippers")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
